$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F, applied to both
# "展览" and "全部类型" worksheets (the two sheets that contain data).
$updates = @{
    2  = 7456
    3  = 7499
    4  = 104
    8  = 119
    11 = 214
    13 = 677
    14 = 594
    16 = 36
    19 = 82
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
